$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark: it used to sit right before the
#    "(4 point) Carilah fungsi alih sistem I!" item; the latest edit
#    happened inside the "Gunakan fungsi bode() ... " item, right
#    after "untuk membuat diagram Bode sistem" (splitting that run in
#    two). Re-adding a bookmark named "_GoBack" automatically removes
#    the previous one (Word only ever keeps a single "_GoBack").
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute("untuk membuat diagram Bode sistem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the Bode-diagram sentence anchor"
}
$splitPos = $anchor.End
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2) Highlight (yellow) the two list items about "sistem Y (gabungan
#    parallel ...)" and "Bandingkan sistem X dan sistem Y ...".
#    These paragraphs previously had no highlight at all; now every
#    run (and the paragraph mark itself) gets wdYellow highlighting,
#    matching the rest of the "Langkah praktikum" list.
# ------------------------------------------------------------------
$paraY = $null
$paraCompare = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*sistem Y (gabungan*") {
        $paraY = $p
    }
    if ($t -like "*Bandingkan sistem X dan sistem Y*") {
        $paraCompare = $p
    }
}

if ($paraY -eq $null) {
    throw "Could not locate the 'sistem Y (gabungan ...)' paragraph"
}
if ($paraCompare -eq $null) {
    throw "Could not locate the 'Bandingkan sistem X dan sistem Y' paragraph"
}

$paraY.Range.Font.HighlightColorIndex = 7
$paraCompare.Range.Font.HighlightColorIndex = 7
